# Update attendance ("想去人数") and min ticket price ("最低票价") figures
# across the 展览, 本地生活 and 全部类型 sheets, reflecting freshly scraped
# numbers (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 -------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 327
$ws.Range("G2").Value = 70
$ws.Range("F5").Value = 175
$ws.Range("F6").Value = 671
$ws.Range("F8").Value = 480
$ws.Range("F9").Value = 85
$ws.Range("F10").Value = 525
$ws.Range("F11").Value = 401
$ws.Range("F14").Value = 116
$ws.Range("F15").Value = 199

# --- Sheet: 本地生活 -----------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6226
$ws.Range("F3").Value = 739
$ws.Range("F5").Value = 1823

# --- Sheet: 全部类型 -----------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6226
$ws.Range("F3").Value = 739
$ws.Range("F5").Value = 1823
$ws.Range("F6").Value = 327
$ws.Range("G6").Value = 70
$ws.Range("F12").Value = 175
$ws.Range("F15").Value = 671
$ws.Range("F19").Value = 480
$ws.Range("F21").Value = 85
$ws.Range("F22").Value = 525
$ws.Range("F24").Value = 401
$ws.Range("F29").Value = 116
$ws.Range("F35").Value = 199
